$d = $word.ActiveDocument
$d.Content.Find.Execute("Nikkelas Raines, Taiki Matehe, Kiarna Broomhead", $false, $false, $false, $false, $false, $true, 1, $false, "Nikkelas Raines, Taiki Matehe", 2)
